$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = "59.702.37"
$ws.Cells.Item(2, 5).Value = "  +1.25%  "

# Row 3
$ws.Cells.Item(3, 4).Value = "2.614.70"
$ws.Cells.Item(3, 5).Value = "  +1.06%  "

# Row 4
$ws.Cells.Item(4, 4).Value = "'1.00"
$ws.Cells.Item(4, 4).Style = "Normal"
$ws.Cells.Item(4, 5).Value = "  +0.16%  "

# Row 5
$ws.Cells.Item(5, 4).Value = "'537.43"
$ws.Cells.Item(5, 4).Style = "Normal"
$ws.Cells.Item(5, 5).Value = "  +3.00%  "

# Row 6
$ws.Cells.Item(6, 4).Value = "'142.29"
$ws.Cells.Item(6, 4).Style = "Normal"
$ws.Cells.Item(6, 5).Value = "  +1.80%  "

# Row 7
$ws.Cells.Item(7, 4).Value = "'1.00"
$ws.Cells.Item(7, 4).Style = "Normal"
$ws.Cells.Item(7, 5).Value = "  +0.18%  "

# Row 8
$ws.Cells.Item(8, 5).Value = "  +0.68%  "

# Row 9
$ws.Cells.Item(9, 5).Value = "  -1.18%  "

# Row 10
$ws.Cells.Item(10, 5).Value = "  +1.95%  "

# Row 11
$ws.Cells.Item(11, 5).Value = "  +1.93%  "

# Row 12
$ws.Cells.Item(12, 4).Value = "'0.133"
$ws.Cells.Item(12, 4).Style = "Normal"
$ws.Cells.Item(12, 5).Value = "  -0.52%  "

# Row 13
$ws.Cells.Item(13, 4).Value = "3.071.51"
$ws.Cells.Item(13, 5).Value = "  +0.92%  "

# Row 14
$ws.Cells.Item(14, 4).Value = "59.666.97"
$ws.Cells.Item(14, 5).Value = "  +1.23%  "

# Row 15
$ws.Cells.Item(15, 4).Value = "'20.78"
$ws.Cells.Item(15, 4).Style = "Normal"
$ws.Cells.Item(15, 5).Value = "  +1.69%  "

# Row 16
$ws.Cells.Item(16, 2).Value = "WrappedEther"
$ws.Cells.Item(16, 3).Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Cells.Item(16, 4).Value = "2.656.01"
$ws.Cells.Item(16, 5).Value = "  +2.83%  "

# Row 17
$ws.Cells.Item(17, 2).Value = "ShibaInu"
$ws.Cells.Item(17, 3).Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Cells.Item(17, 4).Value = "'0.0000134"
$ws.Cells.Item(17, 4).Style = "Normal"
$ws.Cells.Item(17, 5).Value = "  +0.85%  "

# Row 18
$ws.Cells.Item(18, 4).Value = "'343.58"
$ws.Cells.Item(18, 4).Style = "Normal"
$ws.Cells.Item(18, 5).Value = "  +1.51%  "

# Row 19
$ws.Cells.Item(19, 5).Value = "  +1.93%  "

# Row 20
$ws.Cells.Item(20, 5).Value = "  +0.73%  "

# Row 21
$ws.Cells.Item(21, 5).Value = "  -1.42%  "

# Row 22
$ws.Cells.Item(22, 4).Value = "'0.999"
$ws.Cells.Item(22, 4).Style = "Normal"
$ws.Cells.Item(22, 5).Value = "  -0.01%  "

# Row 23
$ws.Cells.Item(23, 4).Value = "'67.60"
$ws.Cells.Item(23, 4).Style = "Normal"
$ws.Cells.Item(23, 5).Value = "  +2.47%  "

# Row 24
$ws.Cells.Item(24, 5).Value = "  +1.75%  "

# Row 25
$ws.Cells.Item(25, 4).Value = "'0.166"
$ws.Cells.Item(25, 4).Style = "Normal"
$ws.Cells.Item(25, 5).Value = "  -0.93%  "

# Row 26
$ws.Cells.Item(26, 4).Value = "'0.997"
$ws.Cells.Item(26, 4).Style = "Normal"
$ws.Cells.Item(26, 5).Value = "  +0.18%  "

# Row 27
$ws.Cells.Item(27, 4).Value = "'7.28"
$ws.Cells.Item(27, 4).Style = "Normal"
$ws.Cells.Item(27, 5).Value = "  +3.73%  "

# Row 28
$ws.Cells.Item(28, 4).Value = "0.0₃0752"
$ws.Cells.Item(28, 5).Value = "  +3.80%  "

# Row 29
$ws.Cells.Item(29, 5).Value = "  +0.02%  "

# Row 30
$ws.Cells.Item(30, 5).Value = "  +6.25%  "

# Row 31
$ws.Cells.Item(31, 4).Value = "'5.86"
$ws.Cells.Item(31, 4).Style = "Normal"
$ws.Cells.Item(31, 5).Value = "  -1.52%  "

# Row 32
$ws.Cells.Item(32, 5).Value = "  +1.47%  "

# Row 33
$ws.Cells.Item(33, 4).Value = "'150.14"
$ws.Cells.Item(33, 4).Style = "Normal"
$ws.Cells.Item(33, 5).Value = "  +0.76%  "

# Row 34
$ws.Cells.Item(34, 5).Value = "  +0.18%  "

# Row 35
$ws.Cells.Item(35, 5).Value = "  +0.22%  "

# Row 36
$ws.Cells.Item(36, 4).Value = "'0.846"
$ws.Cells.Item(36, 4).Style = "Normal"
$ws.Cells.Item(36, 5).Value = "  +4.18%  "

# Row 37
$ws.Cells.Item(37, 5).Value = "  +0.60%  "

# Row 38
$ws.Cells.Item(38, 4).Value = "'0.836"
$ws.Cells.Item(38, 4).Style = "Normal"
$ws.Cells.Item(38, 5).Value = "  +1.04%  "

# Row 39
$ws.Cells.Item(39, 5).Value = "  +1.61%  "

# Row 40
$ws.Cells.Item(40, 4).Value = "'1.00"
$ws.Cells.Item(40, 4).Style = "Normal"
$ws.Cells.Item(40, 5).Value = "  +0.32%  "

# Row 41
$ws.Cells.Item(41, 4).Value = "'276.81"
$ws.Cells.Item(41, 4).Style = "Normal"
$ws.Cells.Item(41, 5).Value = "  +1.51%  "

# Row 42
$ws.Cells.Item(42, 4).Value = "'0.601"
$ws.Cells.Item(42, 4).Style = "Normal"
$ws.Cells.Item(42, 5).Value = "  +1.96%  "

# Row 43
$ws.Cells.Item(43, 2).Value = "WhiteBITCoin"
$ws.Cells.Item(43, 3).Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Cells.Item(43, 4).Value = "'10.75"
$ws.Cells.Item(43, 4).Style = "Normal"
$ws.Cells.Item(43, 5).Value = "  +0.02%  "

# Row 44
$ws.Cells.Item(44, 2).Value = "Stellar"
$ws.Cells.Item(44, 3).Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Cells.Item(44, 4).Value = "'0.0959"
$ws.Cells.Item(44, 4).Style = "Normal"
$ws.Cells.Item(44, 5).Value = "  +0.80%  "

# Row 45
$ws.Cells.Item(45, 4).Value = "'0.0526"
$ws.Cells.Item(45, 4).Style = "Normal"
$ws.Cells.Item(45, 5).Value = "  +1.55%  "

# Row 46
$ws.Cells.Item(46, 4).Value = "1.959.49"
$ws.Cells.Item(46, 5).Value = "  -0.58%  "

# Row 47
$ws.Cells.Item(47, 5).Value = "  +3.34%  "

# Row 48
$ws.Cells.Item(48, 5).Value = "  +2.12%  "

# Row 49
$ws.Cells.Item(49, 5).Value = "  -0.11%  "

# Row 50
$ws.Cells.Item(50, 4).Value = "'111.80"
$ws.Cells.Item(50, 4).Style = "Normal"
$ws.Cells.Item(50, 5).Value = "  -1.44%  "

# Row 51
$ws.Cells.Item(51, 4).Value = "'4.75"
$ws.Cells.Item(51, 4).Style = "Normal"
$ws.Cells.Item(51, 5).Value = "  +0.55%  "
